# Update "想去人数" (column F) values on both the "展览" and "全部类型"
# worksheets. The two sheets contain duplicate data, so the same set of
# row/value updates is applied to each.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    2  = 1938
    4  = 117
    7  = 1612
    9  = 638
    13 = 94
    14 = 223
    17 = 107
    19 = 3737
    21 = 14
    22 = 430
    23 = 344
    24 = 601
    25 = 400
    26 = 351
    28 = 1531
    29 = 11
    30 = 148
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
